$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column J (last existing year column) into column K
# so the new 2023 column matches the style of the preceding years.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new 2023 data column
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 632.9
$ws.Range("K5").Value = 431.8
$ws.Range("K6").Value = 770.7

$excel.CutCopyMode = 0
